# Update Release-Notes.xlsx:
# A new folder entry ("Automated Machine Learning Using AML") was updated
# most recently, so it becomes the new top row of the "Folder Inventory"
# sheet. All existing rows shift down by one, and the stale duplicate
# entry for the same folder (previously the oldest row in the table) is
# removed so the table keeps the same number of data rows.

$wb = $excel.ActiveWorkbook

# --- Folder Inventory sheet -------------------------------------------------
$inv = $wb.Worksheets.Item("Folder Inventory")

# Insert a new row right under the header and populate it with the
# most-recently-updated folder. Clear the formatting Excel copies down
# from the header row on insert, so the new row matches the other
# (unstyled) data rows.
$inv.Rows.Item(2).Insert()
$inv.Rows.Item(2).ClearFormats()
$inv.Range("A2").Value = "Automated Machine Learning Using AML"
$inv.Range("B2").Value = "Automated Machine Learning Using AML"
$inv.Range("C2").Value = "2025-06-12 21:49:46 +0530"
$inv.Range("D2").Value = 1
$inv.Range("E2").Value = "Root"

# The previous (older) row for this same folder now lives at row 46 after
# the insert above shifted everything down by one - remove it so the
# total row count is unchanged.
$inv.Rows.Item(46).Delete()

# --- Metadata sheet ----------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2025-06-12 16:20:02 UTC"
# "Workflow Run" is stored as text in the workbook (not a number), so force
# text entry with a leading apostrophe the same way a user typing into the
# cell would, instead of letting Excel auto-coerce "6" to a number. Then
# reset the cell style back to the plain/unstyled look the rest of the
# sheet uses (the apostrophe entry alone would pick up a "quote prefix"
# style), since the apostrophe itself must not end up in the value.
$meta.Range("B5").Value = "'6"
$meta.Range("B5").Style = $meta.Range("B4").Style

# --- Summary sheet -------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = "2025-06-12 21:49:46 +0530"
